$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order item rows to append (Order_Item_ID, Order_ID, Product_SKU, Quantity, Price, IsDeleted)
$newRows = @(
    @(11, 8,  "CL-G-008",  6,  15000, 0),
    @(12, 9,  "CL-B-006",  10, 10000, 0),
    @(13, 10, "BP-XL-004", 2,  9000,  0),
    @(14, 11, "BP-XL-004", 6,  27000, 0),
    @(15, 12, "BP-L-003",  3,  10500, 0),
    @(16, 13, "BP-M-002",  4,  10000, 0),
    @(17, 14, "BP-L-003",  1,  3500,  0),
    @(18, 15, "BP-L-003",  3,  10500, 0)
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}
